$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footers: the "default" footer (footer2.xml) is Footers.Item(1); the
# "first page" footer (footer1.xml) is Footers.Item(2). Both contain the
# Pearson Edexcel logo inline picture, renamed from image2.png to image1.png.
$ftrDefault = $sec.Footers.Item(1)
if ($ftrDefault.Exists) {
    $ftrDefault.Range.InlineShapes.Item(1).Name = "image1.png"
}
Write-Host "Updated default footer logo name"

$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Exists) {
    $ftrFirst.Range.InlineShapes.Item(1).Name = "image1.png"
}
Write-Host "Updated first-page footer logo name"

# Headers: the "first page" header (header1.xml) is Headers.Item(2) and
# holds the BTEC logo inline picture, renamed from image1.jpg to image2.jpg.
$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Exists) {
    $hdrFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"
}
Write-Host "Updated first-page header logo name"
